# Insert a new weekly price record as row 250 (Acelga, Macroferia Regional de
# Talca) pushing the existing rows 250-268 down to 251-269, matching the
# sheet's "Fruta / hortaliza, semanal" update pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 250..268 down one slot to make room for the new record.
$ws.Rows.Item(250).Insert()

# Seed the new row with the surrounding record's constant fields (Mercado ID,
# Mercado, Region, Codreg, Categoria ID, Categoria, Variedad, Calidad, Unidad
# de comercializacion, Origen, Kg o Unidades, Clasificacion), then overwrite
# the week-specific figures.
$ws.Range("A249:R249").Copy()
$ws.Range("A250").PasteSpecial()

$ws.Range("D250").Value2 = 44746
$ws.Range("J250").Value2 = 500
$ws.Range("K250").Value2 = 2500
$ws.Range("L250").Value2 = 2500
$ws.Range("M250").Value2 = 2500
$ws.Range("P250").Value2 = 625
